$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header - copy header formatting from C1, then center it
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").Value = "PolynomialOrder"
$ws.Range("D1").HorizontalAlignment = -4108   # xlCenter

# Column widths (bestFit-equivalent, nearest achievable grid value)
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 14.33203125

# Updated tStop values and new PolynomialOrder column for each data row
$data = @(
    @{ Row = 2;  B = 0.001;               C = 0.0035;              D = 1 },
    @{ Row = 3;  B = 0.001;               C = 0.003;               D = 1 },
    @{ Row = 4;  B = 0.001;               C = 0.003;               D = 1 },
    @{ Row = 5;  B = 0.001;               C = 0.003;               D = 1 },
    @{ Row = 6;  B = 0.001;               C = 0.003;               D = 1 },
    @{ Row = 7;  B = 0.0011;              C = 0.002;               D = 1 },
    @{ Row = 8;  B = 0.0011;              C = 0.002;               D = 1 },
    @{ Row = 9;  B = 0.001;               C = 0.002;               D = 1 },
    @{ Row = 10; B = 0.001;               C = 0.002;               D = 1 },
    @{ Row = 11; B = 0.0012;              C = 0.002;               D = 1 },
    @{ Row = 12; B = 0.0012;              C = 0.002;               D = 1 },
    @{ Row = 13; B = 0.0012;              C = 0.002;               D = 1 },
    @{ Row = 14; B = 0.0012;              C = 0.002;               D = 1 },
    @{ Row = 15; B = 0.001;               C = 0.0035;              D = 3 },
    @{ Row = 16; B = 0.001;               C = 0.0035;              D = 3 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $dcell = $ws.Cells.Item($r, 4)
    $dcell.Value = $item.D
    $dcell.HorizontalAlignment = -4108   # xlCenter
}

$ws.Range("G11").Select()
